$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph "Added support for RC2014 CP/M format" (the
# last bullet under 1.2c in the Change Log).
# ------------------------------------------------------------------
$findRng = $d.Range(0, $d.Content.End)
$ok = $findRng.Find.Execute("Added support for RC2014 CP/M format")
if (-not $ok) {
    throw "Could not find target paragraph text"
}

$para = $findRng.Paragraphs(1)
$pRng = $para.Range

# Capture the paragraph's own WordOpenXML so we keep its original
# w14:paraId / rsid / pPr attributes intact, then splice in the new
# run layout (split "format" into its own run wrapped in proofErr
# markers) plus the six brand-new paragraphs that follow it in the
# change log.
$origXml = $pRng.WordOpenXML

$oldRunXml = '<w:r><w:t>Added support for RC2014 CP/M format</w:t></w:r>'
$newRunXml = '<w:r><w:t xml:space="preserve">Added support for RC2014 CP/M </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>format</w:t></w:r><w:proofErr w:type="gramEnd"/>'

if ($origXml.IndexOf($oldRunXml) -lt 0) {
    throw "Could not locate expected run XML inside captured paragraph XML"
}

$extraParasXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>1.2c1</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Mostly fixed issue with CP/M disk definition table caused by 1.2c</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>1.2c2</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Finally fixed CP/M </w:t></w:r><w:r><w:t>disk definition table caused by 1.2c</w:t></w:r><w:r><w:t xml:space="preserve"> (mismatch on disk creation button)</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:ind w:left="1080"/></w:pPr></w:p>'

$closeTag = '</w:p>'
$insertPos = $origXml.IndexOf($oldRunXml) + $oldRunXml.Length
$insertPos = $origXml.IndexOf($closeTag, $insertPos) + $closeTag.Length

$newXml = $origXml.Substring(0, $origXml.IndexOf($oldRunXml)) + $newRunXml + $origXml.Substring($origXml.IndexOf($oldRunXml) + $oldRunXml.Length, ($insertPos - ($origXml.IndexOf($oldRunXml) + $oldRunXml.Length))) + $extraParasXml + $origXml.Substring($insertPos)

$pRng.Text = ""
$pRng.InsertXML($newXml)

Write-Output "Applied RC2014/1.2c1/1.2c2 change-log update"
